$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark row 16 ("Binary Tree - getHeight") as reviewed: A16 = "x"
$ws.Cells.Item(16, 1).Value = "x"

# Update the note in C16 with the expanded explanation (two lines)
$newLine = [char]10
$noteText = "每个node都要调用一次getHeight(node) recursive function, 总共有n个node，所以TC = O(n)" + $newLine + "SC = O(height) - blanaced vs ~linkedList"
$ws.Cells.Item(16, 3).Value = $noteText

# Widen column C slightly to fit the new note text
$ws.Columns.Item(3).ColumnWidth = 35.3

# Grow row 16 to fit the new two-line note
$ws.Rows.Item(16).RowHeight = 68

# Update the selection in the bottom-right (frozen) pane to C11
$ws.Range("C11").Select()
